# Remove the "Year of Treatment" column (column B), shifting the
# remaining data columns (old C:K -> new B:J) one position to the left,
# then update the header row text by appending ".deja.deja.deja" to each
# of the shifted header labels (B1:J1). "Country" (A1) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B entirely (this shifts C:K -> B:J and updates the
# sheet dimension from A1:K36 to A1:J36 automatically).
$ws.Columns("B").Delete()

# Append ".deja.deja.deja" to each header cell from B1 to J1 (the
# labels that used to live in C1:K1 before the column shift).
for ($col = 2; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value2 + ".deja.deja.deja"
}
